$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry values for row 71 (append after the last existing row, 70)
$ws.Range("A71").Value = "2025-08-29 06:46:41 UTC"
$ws.Range("B71").Value = "2025-08-29 12:16:41 IST"
$ws.Range("C71").Value = "SKIPPED"
$ws.Range("D71").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E71").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"
$ws.Range("F71").Value = ""
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = ""

# Copy the formatting of the previous row (70) onto the new row (71)
# so the new row matches the existing log styling.
$ws.Range("A70:H70").Copy()
$ws.Range("A71:H71").PasteSpecial(-4122)
